$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing shared strings (A1, B1) and add new ones (C1, D1)
$ws.Range("A1").Value = "ka2xxyys"
$ws.Range("B1").Value = "12g23"
$ws.Range("C1").Value = "karsd23fg"
$ws.Range("D1").Value = "sas234"

# Column widths (target OOXML widths: 20.56, 18.61, 21.39, 20.42 -
# ColumnWidth is stored internally with an implicit +5/6 padding offset,
# so subtract that padding before assigning)
$ws.Columns.Item(1).ColumnWidth = 20.56 - (5/6)
$ws.Columns.Item(2).ColumnWidth = 18.61 - (5/6)
$ws.Columns.Item(3).ColumnWidth = 21.39 - (5/6)
$ws.Columns.Item(4).ColumnWidth = 20.42 - (5/6)

# Selection
$ws.Range("E3").Select()
